$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 3235.3333
$ws.Range("I13").Value = 1177.5
$ws.Range("K13").Value = 1177.5
$ws.Range("M13").Value = -1008.5

$ws.Range("H16").Value = 4499.5
$ws.Range("I16").Value = 4499.5
$ws.Range("K16").Value = 4499.5
$ws.Range("M16").Value = -4269.5

$ws.Range("H19").Value = 861.12
$ws.Range("I19").Value = 530.61536
$ws.Range("J19").Value = 1219.1666
$ws.Range("K19").Value = 530.61536
$ws.Range("L19").Value = 1219.1666
$ws.Range("M19").Value = -355.61536
$ws.Range("N19").Value = -1569.1666

$ws.Range("H33").Value = 3846570.8
$ws.Range("I33").Value = 4167074
$ws.Range("K33").Value = 4167074
$ws.Range("M33").Value = -4166845

$ws.Range("H130").Value = 130745
$ws.Range("J130").Value = 130745
$ws.Range("L130").Value = 130745
$ws.Range("N130").Value = -140785

$ws.Range("H131").Value = 15142.308
$ws.Range("I131").Value = 3349.5454
$ws.Range("K131").Value = 10048.6362
$ws.Range("M131").Value = -5008.636200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 9662.333000000001
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 9662.333000000001
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 9662.333000000001
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -10362.333

$ws.Range("H27").Value = 7497
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 7497
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 7497
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -7865

$ws.Range("H32").Value = 5163.55
$ws.Range("I32").Value = 4640.222
$ws.Range("K32").Value = 4640.222
$ws.Range("M32").Value = -4353.222

$ws.Range("H34").Value = 89995
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 89995
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 89995
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -90537

$ws.Range("H61").Value = 2894.878
$ws.Range("I61").Value = 2443.5312
$ws.Range("J61").Value = 4499.6665
$ws.Range("K61").Value = 2443.5312
$ws.Range("L61").Value = 4499.6665
$ws.Range("M61").Value = -2231.5312
$ws.Range("N61").Value = -4923.6665

$ws.Range("H63").Value = 2614.2856
$ws.Range("I63").Value = 2716.6667
$ws.Range("K63").Value = 2716.6667
$ws.Range("M63").Value = -2030.6667

$ws.Range("H66").Value = 2614.2856
$ws.Range("I66").Value = 2716.6667
$ws.Range("K66").Value = 13583.3335
$ws.Range("M66").Value = -10151.3335

$ws.Range("H74").Value = 11781.667
$ws.Range("I74").Value = 1573.25
$ws.Range("K74").Value = 1573.25
$ws.Range("M74").Value = -699.25

$ws.Range("H77").Value = 11781.667
$ws.Range("I77").Value = 1573.25
$ws.Range("K77").Value = 7866.25
$ws.Range("M77").Value = -3498.25

$ws.Range("H132").Value = 3083.0833
$ws.Range("I132").Value = 2976.8635
$ws.Range("K132").Value = 8930.5905
$ws.Range("M132").Value = -6400.5905

$ws.Range("H136").Value = 2894.878
$ws.Range("I136").Value = 2443.5312
$ws.Range("J136").Value = 4499.6665
$ws.Range("K136").Value = 7330.5936
$ws.Range("L136").Value = 13498.9995
$ws.Range("M136").Value = -4780.5936
$ws.Range("N136").Value = -18598.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 37475
$ws.Range("J63").Value = 37475
$ws.Range("L63").Value = 37475
$ws.Range("N63").Value = -38847

$ws.Range("H66").Value = 37475
$ws.Range("J66").Value = 37475
$ws.Range("L66").Value = 112425
$ws.Range("N66").Value = -119289

$ws.Range("H86").Value = 2035.2941
$ws.Range("I86").Value = 2035.2941
$ws.Range("K86").Value = 2035.2941
$ws.Range("M86").Value = -912.2941000000001

$ws.Range("H89").Value = 2035.2941
$ws.Range("I89").Value = 2035.2941
$ws.Range("K89").Value = 10176.4705
$ws.Range("M89").Value = -4560.470499999999

$ws.Range("H94").Value = 1202
$ws.Range("I94").Value = 1503.4
$ws.Range("K94").Value = 1503.4
$ws.Range("M94").Value = -1052.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 3719.0908
$ws.Range("J15").Value = 4713
$ws.Range("L15").Value = 4713
$ws.Range("N15").Value = -5053

$ws.Range("H29").Value = 25999.4
$ws.Range("J29").Value = 25999.4
$ws.Range("L29").Value = 25999.4
$ws.Range("N29").Value = -26585.4

$ws.Range("H31").Value = 32209.383
$ws.Range("I31").Value = 47636.137
$ws.Range("J31").Value = 3927
$ws.Range("K31").Value = 47636.137
$ws.Range("L31").Value = 3927
$ws.Range("M31").Value = -47341.137
$ws.Range("N31").Value = -4517

$ws.Range("H34").Value = 32209.383
$ws.Range("I34").Value = 47636.137
$ws.Range("J34").Value = 3927
$ws.Range("K34").Value = 47636.137
$ws.Range("L34").Value = 3927
$ws.Range("M34").Value = -47434.137
$ws.Range("N34").Value = -4331

$ws.Range("H43").Value = 15718.2
$ws.Range("J43").Value = 15718.2
$ws.Range("L43").Value = 15718.2
$ws.Range("N43").Value = -16086.2

$ws.Range("H101").Value = 15718.2
$ws.Range("J101").Value = 15718.2
$ws.Range("L101").Value = 15718.2
$ws.Range("N101").Value = -22208.2

$ws.Range("H122").Value = 2061.875
$ws.Range("I122").Value = 1999.2
$ws.Range("J122").Value = 2166.3333
$ws.Range("K122").Value = 5997.6
$ws.Range("L122").Value = 6498.999899999999
$ws.Range("M122").Value = -3547.6
$ws.Range("N122").Value = -11398.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 74.22727
$ws.Range("I2").Value = 73
$ws.Range("K2").Value = 438
$ws.Range("M2").Value = -325

$ws.Range("H92").Value = 1764.6666
$ws.Range("I92").Value = 394
$ws.Range("K92").Value = 1182
$ws.Range("M92").Value = 66

$ws.Range("H102").Value = 287.66666
$ws.Range("I102").Value = 287.66666
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 862.9999799999999
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 1571.00002
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2271.2856
$ws.Range("I122").Value = 1699.8
$ws.Range("K122").Value = 5099.4
$ws.Range("M122").Value = -2649.4

$ws.Range("H126").Value = 3911
$ws.Range("I126").Value = 3469.5625
$ws.Range("J126").Value = 4499.5835
$ws.Range("K126").Value = 10408.6875
$ws.Range("L126").Value = 13498.7505
$ws.Range("M126").Value = -7938.6875
$ws.Range("N126").Value = -18438.7505

$ws.Range("H132").Value = 2711.4167
$ws.Range("I132").Value = 2503.5757
$ws.Range("K132").Value = 7510.7271
$ws.Range("M132").Value = -4980.7271

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7378
$ws.Range("J7").Value = 5749.25
$ws.Range("L7").Value = 5749.25
$ws.Range("N7").Value = -5973.25

$ws.Range("H17").Value = 1538.3334
$ws.Range("I17").Value = 1115.8334
$ws.Range("J17").Value = 2383.3333
$ws.Range("K17").Value = 1115.8334
$ws.Range("L17").Value = 2383.3333
$ws.Range("M17").Value = -945.8334
$ws.Range("N17").Value = -2723.3333

$ws.Range("H43").Value = 10939.521
$ws.Range("I43").Value = 8995
$ws.Range("J43").Value = 17939.8
$ws.Range("K43").Value = 8995
$ws.Range("L43").Value = 17939.8
$ws.Range("M43").Value = -8802
$ws.Range("N43").Value = -18325.8

$ws.Range("H64").Value = 89998
$ws.Range("J64").Value = 89998
$ws.Range("L64").Value = 89998
$ws.Range("N64").Value = -90448

$ws.Range("H67").Value = 89998
$ws.Range("J67").Value = 89998
$ws.Range("L67").Value = 89998
$ws.Range("N67").Value = -91558

$ws.Range("H87").Value = 136728
$ws.Range("J87").Value = 136728
$ws.Range("L87").Value = 136728
$ws.Range("N87").Value = -138974

$ws.Range("H90").Value = 136728
$ws.Range("J90").Value = 136728
$ws.Range("L90").Value = 410184
$ws.Range("N90").Value = -421416

$ws.Range("H126").Value = 7378
$ws.Range("J126").Value = 5749.25
$ws.Range("L126").Value = 17247.75
$ws.Range("N126").Value = -22187.75

$ws.Range("H132").Value = 4705.5
$ws.Range("I132").Value = 4579.6
$ws.Range("K132").Value = 13738.8
$ws.Range("M132").Value = -11208.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 47997.332
$ws.Range("J63").Value = 47997.332
$ws.Range("L63").Value = 47997.332
$ws.Range("N63").Value = -49245.332

$ws.Range("H66").Value = 47997.332
$ws.Range("J66").Value = 47997.332
$ws.Range("L66").Value = 143991.996
$ws.Range("N66").Value = -150231.996

$ws.Range("H132").Value = 1139.375
$ws.Range("I132").Value = 685.8333
$ws.Range("K132").Value = 2057.4999
$ws.Range("M132").Value = 472.5001000000002

$ws.Range("H133").Value = 39886
$ws.Range("J133").Value = 49777.5
$ws.Range("L133").Value = 49777.5
$ws.Range("N133").Value = -59897.5
